# Apply distractor-analysis / unit-test regenerated values to irt_poly_booklet1.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "summary": update distractor analysis values (cols C-L, rows 2-11) ---
$summary = $wb.Worksheets.Item("summary")
$summary.Range("C2").Value = 680
$summary.Range("D2").Value = 644
$summary.Range("E2").Value = 77.02
$summary.Range("F2").Value = -1.5
$summary.Range("H2").Value = 1
$summary.Range("I2").Value = -0.02
$summary.Range("J2").Value = 0.31
$summary.Range("L2").Value = 1.19

$summary.Range("C3").Value = 680
$summary.Range("D3").Value = 624
$summary.Range("I3").Value = -1.08
$summary.Range("J3").Value = 0.38
$summary.Range("K3").Value = 0.06
$summary.Range("L3").Value = 0.87

$summary.Range("C4").Value = 680
$summary.Range("D4").Value = 575
$summary.Range("E4").Value = 41.74
$summary.Range("F4").Value = 0.44
$summary.Range("G4").Value = 0.09
$summary.Range("I4").Value = 0.67
$summary.Range("J4").Value = 0.28
$summary.Range("L4").Value = 0.95

$summary.Range("C5").Value = 680
$summary.Range("D5").Value = 505
$summary.Range("F5").Value = 0.3
$summary.Range("G5").Value = 0.05
$summary.Range("H5").Value = 0.99
$summary.Range("I5").Value = -0.31
$summary.Range("J5").Value = 0.33
$summary.Range("L5").Value = 0.63

$summary.Range("C6").Value = 680
$summary.Range("D6").Value = 349
$summary.Range("E6").Value = 12.03
$summary.Range("F6").Value = 2.45
$summary.Range("G6").Value = 0.18
$summary.Range("H6").Value = 0.97
$summary.Range("I6").Value = -0.18
$summary.Range("J6").Value = 0.27
$summary.Range("L6").Value = 1.62

$summary.Range("C7").Value = 680
$summary.Range("D7").Value = 657
$summary.Range("E7").Value = 75.95
$summary.Range("F7").Value = -1.42
$summary.Range("H7").Value = 1.03
$summary.Range("I7").Value = 0.51
$summary.Range("J7").Value = 0.26
$summary.Range("K7").Value = 0.04
$summary.Range("L7").Value = 0.96

$summary.Range("C8").Value = 680
$summary.Range("D8").Value = 655
$summary.Range("E8").Value = 72.52
$summary.Range("F8").Value = -1.2
$summary.Range("J8").Value = 0.23
$summary.Range("L8").Value = 0.75

$summary.Range("C9").Value = 680
$summary.Range("D9").Value = 647
$summary.Range("E9").Value = 64.91
$summary.Range("F9").Value = -0.78
$summary.Range("H9").Value = 1
$summary.Range("I9").Value = 0.04
$summary.Range("J9").Value = 0.32
$summary.Range("L9").Value = 1.12

$summary.Range("C10").Value = 680
$summary.Range("D10").Value = 648
$summary.Range("E10").Value = 62.35
$summary.Range("F10").Value = -0.63
$summary.Range("I10").Value = 0.79
$summary.Range("J10").Value = 0.28
$summary.Range("L10").Value = 0.97

$summary.Range("C11").Value = 680
$summary.Range("D11").Value = 648
$summary.Range("F11").Value = -1.09
$summary.Range("I11").Value = -0.91
$summary.Range("L11").Value = 0.83

# --- Sheet "model_fit": update N_valid / Deviance / AIC / BIC / EAPrel / WLErel ---
$modelFit = $wb.Worksheets.Item("model_fit")
$modelFit.Range("B2").Value = 680
$modelFit.Range("D2").Value = 7756
$modelFit.Range("E2").Value = 7786
$modelFit.Range("F2").Value = 7854
$modelFit.Range("G2").Value = 0.622
$modelFit.Range("H2").Value = 0.499

$modelFit.Range("B3").Value = 680
$modelFit.Range("D3").Value = 7730
$modelFit.Range("E3").Value = 7778
$modelFit.Range("F3").Value = 7887
$modelFit.Range("G3").Value = 0.635
$modelFit.Range("H3").Value = 0.508

# --- Sheet "steps": update step parameter text values ---
# Leading apostrophe forces these numeric-looking strings to stay text
# (matching the shared-string <t> entries in the original file) instead of
# being auto-coerced to numbers by COM's type inference. Resetting the
# Style afterwards drops the "quote prefix" cell format flag that the
# apostrophe entry leaves behind, so the cell format stays untouched.
$steps = $wb.Worksheets.Item("steps")
$steps.Range("B2").Value = "0.29 (0.12)"
$steps.Range("C2").Value = "'-0.29"
$steps.Range("C2").Style = "Normal"
$steps.Range("B3").Value = "0.93 (0.132)"
$steps.Range("C3").Value = "'-0.93"
$steps.Range("C3").Style = "Normal"
$steps.Range("B4").Value = "0.95 (0.092)"
$steps.Range("C4").Value = "-1.01 (0.092)"
$steps.Range("D4").Value = "'0.06"
$steps.Range("D4").Style = "Normal"
